$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "sports_club_boys"
$ws.Range("C2").Value = "Squash"

$ws.Range("B3").Value = "sports_club_girls"
$ws.Range("C3").Value = "Squash"

$ws.Range("B4").Value = "sports_club_boys"
$ws.Range("C4").Value = "Squash"

$ws.Range("B5").Value = "sports_club_girls"
$ws.Range("C5").Value = "Squash"
